$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column before D ("I/O address") to make room for "I/O device"
$ws.Columns("D").Insert()

$ws.Range("D1").Value = "I/O device"
$ws.Range("D2").Value = "IAS"
$ws.Range("D3").Value = "SSIM"
$ws.Range("D4").Value = "AA"
